$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.522.22"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "1.813.41"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'308.63"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").Value = "'0.3669"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "'0.8800"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'0.07752"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "'19.36"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("D13").Value = "1.843.92"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'5.296"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "'6.372"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "'86.69"
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'0.000008598"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "26.590.42"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "'14.26"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "'5.015"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'1.981"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "'151.43"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'17.93"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").Value = "'112.98"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").Value = "'4.859"
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("D30").Value = "'0.08700"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'3.035"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").Value = "'4.501"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'0.7323"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.120"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'2.678"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").Value = "'1.006"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "'0.05126"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "'2.889"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("D41").Value = "'6.998"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "'0.1558"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "'8.159"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("D45").Value = "'1.008"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'0.4605"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "'9.941"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("D48").Value = "'101.40"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'1.589"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'64.49"
$ws.Range("E51").Value = "  -1.55%  "
